$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2020" column (column R) to the table, mirroring the formatting
# already used by column Q ("2019") for each row, then fill in the new
# year's data values.
$values = [ordered]@{
    4  = 2020
    5  = 5
    6  = 3.5
    7  = 1.8
    8  = 24.4
    9  = 7.2
    10 = 2.9
    11 = 7.4
    12 = 4
    13 = 3.2
    14 = 3.5
}

foreach ($row in $values.Keys) {
    $srcCell = $ws.Range("Q$row")
    $dstCell = $ws.Range("R$row")
    $srcCell.Copy($dstCell) | Out-Null
    $dstCell.Value = $values[$row]
}

# Match the selection left behind by the edit (new column highlighted).
$ws.Range("R4:R14").Select() | Out-Null
